$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) is stored as literal text in the source data (e.g.
# "1.00", "605.10", "3.550.52" -- note some even contain two dots as a
# thousands separator). Prefixing with a leading apostrophe forces Excel to
# keep the assigned value as text instead of silently re-parsing it as a
# number (which would strip trailing zeros / mis-parse multi-dot values).
$textPrefix = "'"

$ws.Range("D2").Value = $textPrefix + '66.390.58'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = $textPrefix + '3.550.52'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("D5").Value = $textPrefix + '605.10'
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("D6").Value = $textPrefix + '144.61'
$ws.Range("E6").Value = '  +0.99%  '
$ws.Range("D7").Value = $textPrefix + '3.550.40'
$ws.Range("E7").Value = '  +0.67%  '
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("E9").Value = '  +2.86%  '
$ws.Range("E10").Value = '  -0.39%  '
$ws.Range("D11").Value = $textPrefix + '7.91'
$ws.Range("E11").Value = '  -1.66%  '
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = $textPrefix + '4.156.14'
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("E14").Value = '  +0.12%  '
$ws.Range("D15").Value = $textPrefix + '30.02'
$ws.Range("E15").Value = '  +0.08%  '
$ws.Range("D16").Value = $textPrefix + '3.557.32'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = $textPrefix + '66.463.45'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").Value = $textPrefix + '11.56'
$ws.Range("E19").Value = '  +5.96%  '
$ws.Range("D20").Value = $textPrefix + '6.17'
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D21").Value = $textPrefix + '14.84'
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("D22").Value = $textPrefix + '431.09'
$ws.Range("E22").Value = '  +1.34%  '
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("D24").Value = $textPrefix + '79.67'
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("D25").Value = $textPrefix + '3.693.49'
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("D26").Value = $textPrefix + '1.00'
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = $textPrefix + '9.19'
$ws.Range("E28").Value = '  +0.60%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = $textPrefix + '2.50'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("B32").Value = 'RenzoRestakedETH'
$ws.Range("C32").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D32").Value = $textPrefix + '3.548.45'
$ws.Range("E32").Value = '  +0.77%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = $textPrefix + '1.45'
$ws.Range("E33").Value = '  -1.78%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = $textPrefix + '25.35'
$ws.Range("E34").Value = '  +0.27%  '
$ws.Range("D35").Value = $textPrefix + '0.152'
$ws.Range("E35").Value = '  -5.34%  '
$ws.Range("D36").Value = $textPrefix + '7.87'
$ws.Range("E36").Value = '  +0.77%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("D40").Value = $textPrefix + '174.83'
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("D41").Value = $textPrefix + '0.0847'
$ws.Range("E41").Value = '  -1.12%  '
$ws.Range("E42").Value = '  +0.21%  '
$ws.Range("D43").Value = $textPrefix + '0.887'
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("E44").Value = '  +1.58%  '
$ws.Range("D45").Value = $textPrefix + '46.05'
$ws.Range("E45").Value = '  +1.46%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("E47").Value = '  +5.52%  '
$ws.Range("D48").Value = $textPrefix + '1.19'
$ws.Range("E48").Value = '  -1.87%  '
$ws.Range("D49").Value = $textPrefix + '25.01'
$ws.Range("E49").Value = '  -3.65%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").Value = $textPrefix + '23.45'
$ws.Range("E51").Value = '  +3.79%  '
